$d = $word.ActiveDocument

$d.Content.Find.Execute("20×87=1740", $true, $false, $false, $false, $false, $true, 1, $false, "24×85=2040", 2) | Out-Null
$d.Content.Find.Execute("27×48=1296", $true, $false, $false, $false, $false, $true, 1, $false, "13×83=1079", 2) | Out-Null
$d.Content.Find.Execute("75×41=3075", $true, $false, $false, $false, $false, $true, 1, $false, "95×33=3135", 2) | Out-Null
$d.Content.Find.Execute("31×78=2418", $true, $false, $false, $false, $false, $true, 1, $false, "89×81=7209", 2) | Out-Null
$d.Content.Find.Execute("91×49=4459", $true, $false, $false, $false, $false, $true, 1, $false, "83×93=7719", 2) | Out-Null
$d.Content.Find.Execute("95×97=9215", $true, $false, $false, $false, $false, $true, 1, $false, "63×37=2331", 2) | Out-Null
$d.Content.Find.Execute("31×85=2635", $true, $false, $false, $false, $false, $true, 1, $false, "56×38=2128", 2) | Out-Null
$d.Content.Find.Execute("42×45=1890", $true, $false, $false, $false, $false, $true, 1, $false, "95×94=8930", 2) | Out-Null
$d.Content.Find.Execute("76×31=2356", $true, $false, $false, $false, $false, $true, 1, $false, "25×97=2425", 2) | Out-Null
$d.Content.Find.Execute("31×36=1116", $true, $false, $false, $false, $false, $true, 1, $false, "86×83=7138", 2) | Out-Null
$d.Content.Find.Execute("37×32=1184", $true, $false, $false, $false, $false, $true, 1, $false, "98×95=9310", 2) | Out-Null
$d.Content.Find.Execute("27×69=1863", $true, $false, $false, $false, $false, $true, 1, $false, "26×97=2522", 2) | Out-Null
$d.Content.Find.Execute("76×78=5928", $true, $false, $false, $false, $false, $true, 1, $false, "63×60=3780", 2) | Out-Null
$d.Content.Find.Execute("34×90=3060", $true, $false, $false, $false, $false, $true, 1, $false, "11×24=264", 2) | Out-Null
$d.Content.Find.Execute("12×47=564", $true, $false, $false, $false, $false, $true, 1, $false, "47×21=987", 2) | Out-Null
$d.Content.Find.Execute("82×81=6642", $true, $false, $false, $false, $false, $true, 1, $false, "65×26=1690", 2) | Out-Null
$d.Content.Find.Execute("72×20=1440", $true, $false, $false, $false, $false, $true, 1, $false, "41×20=820", 2) | Out-Null
$d.Content.Find.Execute("76×24=1824", $true, $false, $false, $false, $false, $true, 1, $false, "93×24=2232", 2) | Out-Null
$d.Content.Find.Execute("71×53=3763", $true, $false, $false, $false, $false, $true, 1, $false, "63×26=1638", 2) | Out-Null
$d.Content.Find.Execute("80×75=6000", $true, $false, $false, $false, $false, $true, 1, $false, "79×60=4740", 2) | Out-Null
$d.Content.Find.Execute("82×91=7462", $true, $false, $false, $false, $false, $true, 1, $false, "32×13=416", 2) | Out-Null
$d.Content.Find.Execute("41×14=574", $true, $false, $false, $false, $false, $true, 1, $false, "39×14=546", 2) | Out-Null
$d.Content.Find.Execute("24×23=552", $true, $false, $false, $false, $false, $true, 1, $false, "64×79=5056", 2) | Out-Null
$d.Content.Find.Execute("19×67=1273", $true, $false, $false, $false, $false, $true, 1, $false, "56×19=1064", 2) | Out-Null
$d.Content.Find.Execute("79×76=6004", $true, $false, $false, $false, $false, $true, 1, $false, "54×89=4806", 2) | Out-Null
